# Update Receptor average/total expression values and their derived
# specificity measures to reflect newly computed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sost -> Lrp5, FAPs -> FAPs): receptor expression values and all
# dependent specificity / edge-weight columns changed because the
# underlying receptor average expression (M2) was recomputed with new TPM.
$ws.Range("M2").Value = 15.01856033333333
$ws.Range("N2").Value = 45.055681
$ws.Range("O2").Value = 0.4908713633047416
$ws.Range("P2").Value = 0.4908713633047417
$ws.Range("Q2").Value = 1.994810239154333
$ws.Range("R2").Value = 17.953292152389
$ws.Range("S2").Value = 0.4908713633047416
$ws.Range("T2").Value = 0.4908713633047417

# Row 3 (Sost -> Lrp5, FAPs -> FAPs): receptor values unchanged but the
# specificity columns shift because row 2's value changed the column sum.
$ws.Range("O3").Value = 0.3099803572711625
$ws.Range("P3").Value = 0.3099803572711625
$ws.Range("S3").Value = 0.3099803572711625
$ws.Range("T3").Value = 0.3099803572711625

# Row 4 (Sost -> Lrp5, FAPs -> MuSCs): same as row 3, specificity-only shift.
$ws.Range("O4").Value = 0.1991482794240958
$ws.Range("P4").Value = 0.1991482794240958
$ws.Range("S4").Value = 0.1991482794240958
$ws.Range("T4").Value = 0.1991482794240958
